# Update res_bus/vm_pu.xlsx results for the "380 kV" case: bus 1 (column B)
# voltage setpoint changed from 1.05 to 1.02 p.u., which changes the
# power-flow results for all buses (columns C-M) across every timestep
# (rows 2-25). Columns A, G, H, N are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.062088900924397
$ws.Range("D2").Value = 1.066775190774129
$ws.Range("E2").Value = 1.057617229651659
$ws.Range("F2").Value = 1.075572017937335
$ws.Range("I2").Value = 1.052825072051579
$ws.Range("J2").Value = 1.067061224130507
$ws.Range("K2").Value = 1.069484833642416
$ws.Range("L2").Value = 1.060351709029094
$ws.Range("M2").Value = 1.078258234190973

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.063317741591879
$ws.Range("D3").Value = 1.067772924698103
$ws.Range("E3").Value = 1.05867543085846
$ws.Range("F3").Value = 1.076701715084911
$ws.Range("I3").Value = 1.053222933430652
$ws.Range("J3").Value = 1.067943031481683
$ws.Range("K3").Value = 1.070297771097248
$ws.Range("L3").Value = 1.06122316132055
$ws.Range("M3").Value = 1.079204506370166

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.064112700855984
$ws.Range("D4").Value = 1.06841832831858
$ws.Range("E4").Value = 1.059360216827389
$ws.Range("F4").Value = 1.077432804832294
$ws.Range("I4").Value = 1.053479079150611
$ws.Range("J4").Value = 1.06851289274186
$ws.Range("K4").Value = 1.070822983240496
$ws.Range("L4").Value = 1.06178650676341
$ws.Range("M4").Value = 1.079816316361867

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.064446859835625
$ws.Range("D5").Value = 1.068689609687074
$ws.Range("E5").Value = 1.059648115906482
$ws.Range("F5").Value = 1.077740179960075
$ws.Range("I5").Value = 1.053586452756275
$ws.Range("J5").Value = 1.068752289398869
$ws.Range("K5").Value = 1.071043588724869
$ws.Range("L5").Value = 1.062023208226536
$ws.Range("M5").Value = 1.080073404597423

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.064502964196398
$ws.Range("D6").Value = 1.068735156330973
$ws.Range("E6").Value = 1.0596964563386
$ws.Range("F6").Value = 1.077791791092139
$ws.Range("I6").Value = 1.05360446308981
$ws.Range("J6").Value = 1.068792475035522
$ws.Range("K6").Value = 1.071080618023098
$ws.Range("L6").Value = 1.062062943907718
$ws.Range("M6").Value = 1.080116564028673

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.064117166066895
$ws.Range("D7").Value = 1.068421953374439
$ws.Range("E7").Value = 1.059364063688649
$ws.Range("F7").Value = 1.07743691189429
$ws.Range("I7").Value = 1.053480515099667
$ws.Range("J7").Value = 1.068516092250242
$ws.Range("K7").Value = 1.070825931741795
$ws.Range("L7").Value = 1.061789670085544
$ws.Range("M7").Value = 1.079819752044753

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.062504232361558
$ws.Range("D8").Value = 1.067112420907356
$ws.Range("E8").Value = 1.057974842039847
$ws.Range("F8").Value = 1.07595378425559
$ws.Range("I8").Value = 1.05295980005934
$ws.Range("J8").Value = 1.067359386224279
$ws.Range("K8").Value = 1.069759738472581
$ws.Range("L8").Value = 1.060646333276851
$ws.Range("M8").Value = 1.078578133622756

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.059660552935182
$ws.Range("D9").Value = 1.064803311159054
$ws.Range("E9").Value = 1.055527266182552
$ws.Range("F9").Value = 1.073341044753619
$ws.Range("I9").Value = 1.052032276235903
$ws.Range("J9").Value = 1.065315506977774
$ws.Range("K9").Value = 1.067874705571736
$ws.Range("L9").Value = 1.058627432233179
$ws.Range("M9").Value = 1.076386436853937

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.057763649128561
$ws.Range("D10").Value = 1.063262809491461
$ws.Range("E10").Value = 1.053895763530794
$ws.Range("F10").Value = 1.07159963585963
$ws.Range("I10").Value = 1.051407199565406
$ws.Range("J10").Value = 1.063949077355302
$ws.Range("K10").Value = 1.066613747973537
$ws.Range("L10").Value = 1.057278617955635
$ws.Range("M10").Value = 1.074922683793628

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056941976444259
$ws.Range("D11").Value = 1.062595482308249
$ws.Range("E11").Value = 1.053189341329666
$ws.Range("F11").Value = 1.07084566977834
$ws.Range("I11").Value = 1.05113493110197
$ws.Range("J11").Value = 1.063356470803053
$ws.Range("K11").Value = 1.06606671340387
$ws.Range("L11").Value = 1.056693869389604
$ws.Range("M11").Value = 1.074288225982126

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.056636722978922
$ws.Range("D12").Value = 1.062347564005918
$ws.Range("E12").Value = 1.052926947717615
$ws.Range("F12").Value = 1.070565623077571
$ws.Range("I12").Value = 1.051033556436815
$ws.Range("J12").Value = 1.063136208636866
$ws.Range("K12").Value = 1.065863364167184
$ws.Range("L12").Value = 1.056476560812769
$ws.Range("M12").Value = 1.074052462121572

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.05670220307889
$ws.Range("D13").Value = 1.062400745301493
$ws.Range("E13").Value = 1.052983231921832
$ws.Range("F13").Value = 1.070625693622359
$ws.Range("I13").Value = 1.051055312610951
$ws.Range("J13").Value = 1.063183462055185
$ws.Range("K13").Value = 1.065906990359118
$ws.Range("L13").Value = 1.056523179099483
$ws.Range("M13").Value = 1.074103038725742

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.056916745083727
$ws.Range("D14").Value = 1.062574990193983
$ws.Range("E14").Value = 1.05316765173746
$ws.Range("F14").Value = 1.07082252083779
$ws.Range("I14").Value = 1.051126556382193
$ws.Range("J14").Value = 1.063338266762903
$ws.Range("K14").Value = 1.066049907687053
$ws.Range("L14").Value = 1.05667590878355
$ws.Range("M14").Value = 1.074268739674536

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.057048925106846
$ws.Range("D15").Value = 1.062682342451776
$ws.Range("E15").Value = 1.053281279236025
$ws.Range("F15").Value = 1.070943793845608
$ws.Range("I15").Value = 1.051170419927139
$ws.Range("J15").Value = 1.063433628231342
$ws.Range("K15").Value = 1.066137943028706
$ws.Range("L15").Value = 1.056769996375488
$ws.Range("M15").Value = 1.074370820473709

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057818174291485
$ws.Range("D16").Value = 1.063307091798305
$ws.Range("E16").Value = 1.053942646946765
$ws.Range("F16").Value = 1.071649675570239
$ws.Range("I16").Value = 1.051425235224744
$ws.Range("J16").Value = 1.063988386908155
$ws.Range("K16").Value = 1.066650031014592
$ws.Range("L16").Value = 1.057317410868674
$ws.Range("M16").Value = 1.074964777003701

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.058300621714147
$ws.Range("D17").Value = 1.063698904775459
$ws.Range("E17").Value = 1.054357512312995
$ws.Range("F17").Value = 1.072092475696355
$ws.Range("I17").Value = 1.05158464360345
$ws.Range("J17").Value = 1.064336121414021
$ws.Range("K17").Value = 1.066970973203601
$ws.Range("L17").Value = 1.057660600369903
$ws.Range("M17").Value = 1.075337177109547

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.0585819962757
$ws.Range("D18").Value = 1.06392741574487
$ws.Range("E18").Value = 1.054599499172869
$ws.Range("F18").Value = 1.072350760899677
$ws.Range("I18").Value = 1.051677468894496
$ws.Range("J18").Value = 1.064538859097683
$ws.Range("K18").Value = 1.067158073903058
$ws.Range("L18").Value = 1.057860709250883
$ws.Range("M18").Value = 1.07555432983664

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.058677932908664
$ws.Range("D19").Value = 1.064005327552498
$ws.Range("E19").Value = 1.054682010972366
$ws.Range("F19").Value = 1.072438830775094
$ws.Range("I19").Value = 1.051709093664788
$ws.Range("J19").Value = 1.064607972191707
$ws.Range("K19").Value = 1.06722185357558
$ws.Range("L19").Value = 1.057928929755321
$ws.Range("M19").Value = 1.075628362837164

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.058248862675352
$ws.Range("D20").Value = 1.06365686974452
$ws.Range("E20").Value = 1.054313000895042
$ws.Range("F20").Value = 1.072044966681572
$ws.Range("I20").Value = 1.051567556622849
$ws.Range("J20").Value = 1.064298822109971
$ws.Range("K20").Value = 1.06693654944051
$ws.Range("L20").Value = 1.057623786402123
$ws.Range("M20").Value = 1.07529722852155

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056853569145567
$ws.Range("D21").Value = 1.062523680620611
$ws.Range("E21").Value = 1.053113344624989
$ws.Range("F21").Value = 1.0707645598679
$ws.Range("I21").Value = 1.05110558355449
$ws.Range("J21").Value = 1.063292684558411
$ws.Range("K21").Value = 1.066007826419186
$ws.Range("L21").Value = 1.056630936660292
$ws.Range("M21").Value = 1.074219947590813

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.055976015847032
$ws.Range("D22").Value = 1.061810948658118
$ws.Range("E22").Value = 1.052359090074913
$ws.Range("F22").Value = 1.069959573674146
$ws.Range("I22").Value = 1.050813722048373
$ws.Range("J22").Value = 1.062659265350227
$ws.Range("K22").Value = 1.065422996811269
$ws.Range("L22").Value = 1.056006073352793
$ws.Range("M22").Value = 1.073542050751624

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.056441250518677
$ws.Range("D23").Value = 1.062188805561643
$ws.Range("E23").Value = 1.052758933524499
$ws.Range("F23").Value = 1.070386307231641
$ws.Range("I23").Value = 1.050968576369558
$ws.Range("J23").Value = 1.062995131127882
$ws.Range("K23").Value = 1.065733112215189
$ws.Range("L23").Value = 1.056337384347201
$ws.Range("M23").Value = 1.073901470947332

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.058272250446238
$ws.Range("D24").Value = 1.063675863648892
$ws.Range("E24").Value = 1.054333113679963
$ws.Range("F24").Value = 1.072066433937662
$ws.Range("I24").Value = 1.051575277972601
$ws.Range("J24").Value = 1.064315676338004
$ws.Range("K24").Value = 1.066952104364657
$ws.Range("L24").Value = 1.057640421258938
$ws.Range("M24").Value = 1.07531527976042

$ws.Range("B25").Value = 1.019999999999999
$ws.Range("C25").Value = 1.060395900423877
$ws.Range("D25").Value = 1.0654004603227
$ws.Range("E25").Value = 1.056159980889576
$ws.Range("F25").Value = 1.074016422024656
$ws.Range("I25").Value = 1.052273246741163
$ws.Range("J25").Value = 1.065844571201166
$ws.Range("K25").Value = 1.068362779988752
$ws.Range("L25").Value = 1.059149869844057
$ws.Range("M25").Value = 1.076953500528012
